# Apply the "456a3b4" gh-pages data refresh to 合肥-漫展信息.xlsx
#
# Updates "想去人数" (want-to-go count) figures on both the "展览"
# (Exhibition) sheet and the "全部类型" (All types) sheet, marks the
# "合肥·WA二次元饭局" event as cancelled / not-for-sale, and appends a
# newly scraped event ("合肥·环形宇宙动漫游戏嘉年华第7届") to both
# sheets.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, $text)
    # Assigning a plain string through .Value lets the engine "helpfully"
    # auto-detect dates/numbers (e.g. "2024-08-03" -> a date serial). Force
    # the cell to text first so the literal string is preserved, then drop
    # the now-unneeded number format so the cell is left unstyled, matching
    # the plain inlineStr cells used elsewhere in this workbook.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Update-CommonRows {
    param($ws)

    # "想去人数" (want-to-go count) bumps shared by both sheets.
    $ws.Cells.Item(3, 6).Value = 418
    $ws.Cells.Item(4, 6).Value = 12199
    $ws.Cells.Item(5, 6).Value = 1262
    $ws.Cells.Item(9, 6).Value = 157
    $ws.Cells.Item(11, 6).Value = 442

    # "合肥·WA二次元饭局" got cancelled.
    $ws.Cells.Item(13, 3).Value = "合肥·WA二次元饭局（取消）"
    $ws.Cells.Item(13, 6).Value = 62
    $ws.Cells.Item(13, 7).Value = "不可售"

    $ws.Cells.Item(15, 6).Value = 38
    $ws.Cells.Item(16, 6).Value = 361
    $ws.Cells.Item(17, 6).Value = 2799
    $ws.Cells.Item(19, 6).Value = 931
    $ws.Cells.Item(20, 6).Value = 122
}

function Add-HuanXingYuZhouRow {
    param($ws, $rowIndex, $indexValue)

    # Match the bold/centered/bordered style used by the rest of column A.
    $ws.Cells.Item($rowIndex - 1, 1).Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($rowIndex, 1).Value = $indexValue
    Set-TextValue $ws.Cells.Item($rowIndex, 2) "2024-08-03"
    $ws.Cells.Item($rowIndex, 3).Value = "合肥·环形宇宙动漫游戏嘉年华第7届"
    $ws.Cells.Item($rowIndex, 4).Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Cells.Item($rowIndex, 5).Value = "2024.08.03 09:30-08.04 17:00"
    $ws.Cells.Item($rowIndex, 6).Value = 0
    $ws.Cells.Item($rowIndex, 7).Value = 70
    $ws.Cells.Item($rowIndex, 8).Value = "https://show.bilibili.com/platform/detail.html?id=84767"
    $ws.Cells.Item($rowIndex, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"
}

# ---- Sheet "展览" (Exhibition) -------------------------------------------
$wsExhibition = $wb.Worksheets.Item("展览")
Update-CommonRows $wsExhibition
Add-HuanXingYuZhouRow $wsExhibition 21 20

# ---- Sheet "全部类型" (All types) -----------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
Update-CommonRows $wsAll

# Insert a new row 21 (pushing the existing row 21 down to row 22, carrying
# its formatting/content with it) and fill it in with the new event.
$wsAll.Rows.Item(21).Insert()
Add-HuanXingYuZhouRow $wsAll 21 20

# The event that used to be row 21 ("合肥·首届包河留声机音乐节...") is now
# row 22; bump its index number to match its new position.
$wsAll.Cells.Item(22, 1).Value = 21
